$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.930.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.155.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.154.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("E10").Value = "  -2.45%  "

$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("E12").Value = "  -0.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.673.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.897.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.159.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "504.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.45%  "

$ws.Range("E21").Value = "  -0.38%  "

$ws.Range("E22").Value = "  -2.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.09%  "

$ws.Range("E24").Value = "  -1.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.34%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("E30").Value = "  +4.07%  "

$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.02%  "

$ws.Range("E35").Value = "  -1.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0891"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "479.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0415"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.19%  "

$ws.Range("E40").Value = "  -2.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.990.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.11%  "

$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("E44").Value = "  -4.19%  "

$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0592"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.39%  "

$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.73%  "

$ws.Range("E51").Value = "  +14.21%  "
